$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Columns that differ between row 13 and row 14 and need to be swapped
$cols = @("A", "B", "E", "F", "G", "H", "P", "Q", "R", "S")

foreach ($col in $cols) {
    $addr13 = "$col`13"
    $addr14 = "$col`14"
    $val13 = $ws.Range($addr13).Value2
    $val14 = $ws.Range($addr14).Value2
    $ws.Range($addr13).Value = $val14
    $ws.Range($addr14).Value = $val13
}
